# Weekly fruit/vegetable price update: insert a new daily price record
# for "Ajo" (Vega Monumental Concepción) as the new most-recent row.
#
# Effect: a brand-new row is inserted at row 40, pushing the previous
# rows 40..162 down to 41..163 (dimension grows from A1:R162 to A1:R163),
# and the new row 40 is populated with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 40 (shifts existing rows down).
$ws.Rows.Item(40).Insert()

# Populate the new row with the latest price record.
$ws.Cells.Item(40, 1).Value  = 11
$ws.Cells.Item(40, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value  = "Bíobío"
$ws.Cells.Item(40, 4).Value  = 44672
$ws.Cells.Item(40, 5).Value  = 8
$ws.Cells.Item(40, 6).Value  = 100112003
$ws.Cells.Item(40, 7).Value  = "Ajo"
$ws.Cells.Item(40, 8).Value  = "Chino"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 270
$ws.Cells.Item(40, 11).Value = 17000
$ws.Cells.Item(40, 12).Value = 18000
$ws.Cells.Item(40, 13).Value = 17556
$ws.Cells.Item(40, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(40, 15).Value = "China"
$ws.Cells.Item(40, 16).Value = 1756
$ws.Cells.Item(40, 17).Value = 10
$ws.Cells.Item(40, 18).Value = "Hortaliza"
